$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.821.53'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.557.24'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.50%  '
$ws.Range('E7').Value = '  -0.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.531'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.63'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0808'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.44'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.952.20'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.34%  '
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.92'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.619.59'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.839'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.825.60'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.76'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0960'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '247.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.05'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.15%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.38'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '39.90'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.85%  '
$ws.Range('E30').Value = '  -1.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '158.29'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.77'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0797'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.08%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.30'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.98%  '
$ws.Range('E35').Value = '  -2.39%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.66'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.83%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.77'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.65%  '
$ws.Range('E38').Value = '  +11.96%  '
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.90'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.09'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.47%  '
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('E44').Value = '  -1.08%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.20'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.25%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.990.78'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.96'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.806.32'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '81.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.30%  '
$ws.Range('E50').Value = '  +2.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.64'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.71%  '
